$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 38: Review af OC0803 / Reviewer ---
$ws.Range("A38").Value = "Review af OC0803"
$ws.Range("B38").Value = "Reviewer"
$ws.Range("C38").Value = 43900
$ws.Range("D38").Value = 0.4375
$ws.Range("E38").Value = 0.46180555555555558
$ws.Range("G38").Formula = "=E38-D38"

# --- Row 39: Systemtest OC0802 / Tester ---
$ws.Range("A39").Value = "Systemtest OC0802"
$ws.Range("B39").Value = "Tester"
$ws.Range("C39").Value = 43900
$ws.Range("D39").Value = 0.4861111111111111
$ws.Range("E39").Value = 0.5
$ws.Range("G39").Formula = "=E39-D39"

# --- Row 40: AEndringer af SSD10 i review med Marc / Software Architect ---
$ws.Range("A40").Value = "Ændringer af SSD10 i review med Marc"
$ws.Range("B40").Value = "Software Architect"
$ws.Range("C40").Value = 43900
$ws.Range("D40").Value = 0.51736111111111105
$ws.Range("E40").Value = 0.65277777777777779
$ws.Range("G40").Formula = "=E40-D40"

# --- Row 41: Opstart af branch til systemtest 0802 med ANDP / Any Role ---
$ws.Range("A41").Value = "Opstart af branch til systemtest 0802 med ANDP"
$ws.Range("B41").Value = "Any Role"
$ws.Range("C41").Value = 43900
$ws.Range("D41").Value = 0.65277777777777779
$ws.Range("E41").Value = 0.65972222222222221
$ws.Range("G41").Formula = "=E41-D41"

# --- Update the view: scroll position and active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B43").Select()

$wb.Save()
